# Auto-generated Word COM-interop script to update date and multiplication table values
$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-02-03 Saturday", $true, $true, $false, $false, $false, $true, 1, $false, "2024-02-04 Sunday", 2) | Out-Null
$d.Content.Find.Execute("95×23=2185", $true, $true, $false, $false, $false, $true, 1, $false, "44×12=528", 2) | Out-Null
$d.Content.Find.Execute("80×73=5840", $true, $true, $false, $false, $false, $true, 1, $false, "37×38=1406", 2) | Out-Null
$d.Content.Find.Execute("68×75=5100", $true, $true, $false, $false, $false, $true, 1, $false, "95×87=8265", 2) | Out-Null
$d.Content.Find.Execute("69×86=5934", $true, $true, $false, $false, $false, $true, 1, $false, "61×56=3416", 2) | Out-Null
$d.Content.Find.Execute("30×73=2190", $true, $true, $false, $false, $false, $true, 1, $false, "90×67=6030", 2) | Out-Null
$d.Content.Find.Execute("18×67=1206", $true, $true, $false, $false, $false, $true, 1, $false, "98×68=6664", 2) | Out-Null
$d.Content.Find.Execute("93×47=4371", $true, $true, $false, $false, $false, $true, 1, $false, "63×42=2646", 2) | Out-Null
$d.Content.Find.Execute("52×62=3224", $true, $true, $false, $false, $false, $true, 1, $false, "40×41=1640", 2) | Out-Null
$d.Content.Find.Execute("47×51=2397", $true, $true, $false, $false, $false, $true, 1, $false, "67×23=1541", 2) | Out-Null
$d.Content.Find.Execute("12×54=648", $true, $true, $false, $false, $false, $true, 1, $false, "64×50=3200", 2) | Out-Null
$d.Content.Find.Execute("47×47=2209", $true, $true, $false, $false, $false, $true, 1, $false, "81×50=4050", 2) | Out-Null
$d.Content.Find.Execute("97×16=1552", $true, $true, $false, $false, $false, $true, 1, $false, "27×35=945", 2) | Out-Null
$d.Content.Find.Execute("92×17=1564", $true, $true, $false, $false, $false, $true, 1, $false, "43×49=2107", 2) | Out-Null
$d.Content.Find.Execute("85×75=6375", $true, $true, $false, $false, $false, $true, 1, $false, "71×99=7029", 2) | Out-Null
$d.Content.Find.Execute("94×25=2350", $true, $true, $false, $false, $false, $true, 1, $false, "62×52=3224", 2) | Out-Null
$d.Content.Find.Execute("97×70=6790", $true, $true, $false, $false, $false, $true, 1, $false, "33×78=2574", 2) | Out-Null
$d.Content.Find.Execute("61×17=1037", $true, $true, $false, $false, $false, $true, 1, $false, "75×59=4425", 2) | Out-Null
$d.Content.Find.Execute("30×33=990", $true, $true, $false, $false, $false, $true, 1, $false, "64×18=1152", 2) | Out-Null
$d.Content.Find.Execute("67×13=871", $true, $true, $false, $false, $false, $true, 1, $false, "93×30=2790", 2) | Out-Null
$d.Content.Find.Execute("21×60=1260", $true, $true, $false, $false, $false, $true, 1, $false, "11×47=517", 2) | Out-Null
$d.Content.Find.Execute("91×18=1638", $true, $true, $false, $false, $false, $true, 1, $false, "50×70=3500", 2) | Out-Null
$d.Content.Find.Execute("75×43=3225", $true, $true, $false, $false, $false, $true, 1, $false, "50×37=1850", 2) | Out-Null
$d.Content.Find.Execute("75×67=5025", $true, $true, $false, $false, $false, $true, 1, $false, "39×42=1638", 2) | Out-Null
$d.Content.Find.Execute("89×79=7031", $true, $true, $false, $false, $false, $true, 1, $false, "28×41=1148", 2) | Out-Null
$d.Content.Find.Execute("54×39=2106", $true, $true, $false, $false, $false, $true, 1, $false, "69×40=2760", 2) | Out-Null
